$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current (before) layout:
#   TC2 block step2 (B20/D20): "ordenar pelo nome do servidor" (ordering)
#   TC3 block step2 (B28/D28): "Indica parametros...Filtra" (filter)
#   TC4 block step2 (B36/D36): "cancelamento de uma diaria" (cancel)
#
# Target (after) layout per diff:
#   TC2 block step2 (B20/D20): filter content
#   TC3 block step2 (B28/D28): cancel content
#   TC4 block step2 (B36/D36): ordering content

$orderingB = $ws.Range("B20").Text
$orderingD = $ws.Range("D20").Text

$filterB = $ws.Range("B28").Text
$filterD = $ws.Range("D28").Text

$cancelB = $ws.Range("B36").Text
$cancelD = $ws.Range("D36").Text

# Move filter content into TC2's step2 slot
$ws.Range("B20").Value = $filterB
$ws.Range("D20").Value = $filterD

# Move cancel content into TC3's step2 slot
$ws.Range("B28").Value = $cancelB
$ws.Range("D28").Value = $cancelD

# Move ordering content into TC4's step2 slot
$ws.Range("B36").Value = $orderingB
$ws.Range("D36").Value = $orderingD
